# MFY auto commit at 01/12/2021 14:06:58
# Adds a new "exam date/time" column (inserted before the old column F,
# which shifts right to become column G), fills it in for each of the
# five subject blocks, adjusts page setup / margins / zoom, and moves
# the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before F; old F (room numbers) shifts to G ---
$ws.Columns("F").Insert()

# --- New column F: exam date / time for each subject block ---
$ws.Range("F3").Value = "Tuesday, 30 Nov 11:30 to 1:00"
$ws.Range("F9").Value = "Firday, 3 Dec 9:30 to 11:00"
$ws.Range("F6").Value = "Wednesday, 1 Dec 11:30 to 1:00"
$ws.Range("F12").Value = "Saturday, 4 Dec 9:30 to 11:00"
$ws.Range("F15").Value = "Sunday, 30 Dec 2:00 to 3:30"

# Match the formatting already used across the table (centered, wrapped)
# by copying it from column A onto the new column F, then merge each
# subject's two-row block just like the other columns.
$ws.Range("A3:A16").Copy()
$ws.Range("F3:F16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F3:F4").Merge()
$ws.Range("F6:F7").Merge()
$ws.Range("F9:F10").Merge()
$ws.Range("F12:F13").Merge()
$ws.Range("F15:F16").Merge()

# New column's width
$ws.Columns("F").ColumnWidth = 18.59

# --- Page setup: narrower margins, fit-to-width printing, smaller paper scale ---
$ws.PageSetup.LeftMargin = 18
$ws.PageSetup.RightMargin = 18
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Zoom = 77
$ws.PageSetup.FitToPagesTall = $false

# --- Sheet view: zoom in and move the selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 115
$ws.Range("B2").Select()
